$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.026.75"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.516.49"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "531.22"
$ws.Range("E5").Value = "  -1.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.33"
$ws.Range("E6").Value = "  -3.06%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  -1.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.520.27"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("E12").Value = "  -2.46%  "
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.961.61"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.05"
$ws.Range("E15").Value = "  -1.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.981.59"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.512.14"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.03"
$ws.Range("E19").Value = "  -1.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.26"
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.25"
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.82"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.25"
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("E25").Value = "  -3.67%  "
$ws.Range("E26").Value = "  +1.81%  "
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.81"
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.84"
$ws.Range("E29").Value = "  +2.70%  "
$ws.Range("E30").Value = "  -0.60%  "
$ws.Range("E31").Value = "  -1.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "161.33"
$ws.Range("E32").Value = "  +3.10%  "
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.13"
$ws.Range("E34").Value = "  -5.56%  "
$ws.Range("E35").Value = "  -1.08%  "
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.22"
$ws.Range("E37").Value = "  -2.90%  "
$ws.Range("E38").Value = "  -1.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.97"
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("E41").Value = "  -2.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.21"
$ws.Range("E42").Value = "  -8.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "280.64"
$ws.Range("E43").Value = "  -5.28%  "
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0932"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.36"
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.42"
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0509"
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0223"
$ws.Range("E51").Value = "  -2.31%  "
